$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.160.04'
$ws.Range("E2").Value = '  +3.28%  '
$ws.Range("D3").Value = '2.621.28'
$ws.Range("E3").Value = '  +3.82%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = "'604.94"
$ws.Range("E5").Value = '  +1.63%  '
$ws.Range("D6").Value = "'179.53"
$ws.Range("E6").Value = '  +1.11%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  +0.73%  '
$ws.Range("D9").Value = '2.619.92'
$ws.Range("E9").Value = '  +3.77%  '
$ws.Range("D10").Value = "'0.166"
$ws.Range("E10").Value = '  +12.60%  '
$ws.Range("E11").Value = '  +0.37%  '
$ws.Range("E12").Value = '  +2.18%  '
$ws.Range("D13").Value = "'5.01"
$ws.Range("E13").Value = '  +0.38%  '
$ws.Range("D14").Value = '3.110.21'
$ws.Range("E14").Value = '  +4.15%  '
$ws.Range("E15").Value = '  +8.65%  '
$ws.Range("D16").Value = "'26.63"
$ws.Range("E16").Value = '  +1.66%  '
$ws.Range("D17").Value = '71.175.55'
$ws.Range("E17").Value = '  +3.40%  '
$ws.Range("D18").Value = '2.594.93'
$ws.Range("E18").Value = '  +2.65%  '
$ws.Range("D19").Value = "'378.78"
$ws.Range("E19").Value = '  +6.68%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").Value = "'11.49"
$ws.Range("E20").Value = '  +3.07%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = "'7.86"
$ws.Range("E21").Value = '  +4.01%  '
$ws.Range("D22").Value = "'4.16"
$ws.Range("E22").Value = '  -0.57%  '
$ws.Range("E23").Value = '  +17.44%  '
$ws.Range("D24").Value = "'72.32"
$ws.Range("E24").Value = '  +2.44%  '
$ws.Range("D25").Value = "'4.47"
$ws.Range("E25").Value = '  +6.00%  '
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("D27").Value = "'9.98"
$ws.Range("E27").Value = '  +10.93%  '
$ws.Range("D28").Value = '2.758.87'
$ws.Range("E28").Value = '  +3.88%  '
$ws.Range("E29").Value = '  +0.25%  '
$ws.Range("D30").Value = "'547.93"
$ws.Range("E30").Value = '  +4.27%  '
$ws.Range("D31").Value = '0.0₃0953'
$ws.Range("E31").Value = '  +6.91%  '
$ws.Range("D32").Value = "'8.06"
$ws.Range("E32").Value = '  +3.48%  '
$ws.Range("D33").Value = "'1.33"
$ws.Range("E33").Value = '  +6.64%  '
$ws.Range("E34").Value = '  +2.43%  '
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = '  -0.05%  '
$ws.Range("D36").Value = "'165.87"
$ws.Range("E36").Value = '  +2.35%  '
$ws.Range("D37").Value = "'0.116"
$ws.Range("E37").Value = '  -3.30%  '
$ws.Range("D38").Value = "'19.21"
$ws.Range("E38").Value = '  +4.27%  '
$ws.Range("B39").Value = 'WhiteBITCoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D39").Value = "'19.16"
$ws.Range("E39").Value = '  +2.48%  '
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").Value = "'1.88"
$ws.Range("E40").Value = '  +5.65%  '
$ws.Range("E41").Value = '  +5.67%  '
$ws.Range("E42").Value = '  -0.06%  '
$ws.Range("E43").Value = '  +9.09%  '
$ws.Range("D44").Value = "'5.03"
$ws.Range("E44").Value = '  +3.70%  '
$ws.Range("D45").Value = "'0.332"
$ws.Range("E45").Value = '  +1.59%  '
$ws.Range("D46").Value = "'39.95"
$ws.Range("E46").Value = '  +2.44%  '
$ws.Range("D47").Value = "'152.95"
$ws.Range("E47").Value = '  -0.16%  '
$ws.Range("E48").Value = '  +1.57%  '
$ws.Range("D49").Value = "'0.534"
$ws.Range("E49").Value = '  +2.75%  '
$ws.Range("E50").Value = '  +4.78%  '
$ws.Range("D51").Value = '0.0₆0261'
$ws.Range("E51").Value = '  +2.79%  '
